# Updated cryptos list data (Price/Volume columns) for rows 2-51.
# Values that look numeric are force-written as text (NumberFormat "@")
# then ClearFormats() restores the default (unstyled) cell formatting so
# the stored style index is unchanged, while the literal text is preserved
# exactly (avoids Excel auto-converting e.g. "0.1140" -> 0.114).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.323.13"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").Value = "1.839.19"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9984"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.23"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6263"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07428"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.89%  "

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.83"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07721"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.10%  "

$ws.Range("D12").Value = "1.832.75"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.954"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6743"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001023"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.78"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.238"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.37%  "

$ws.Range("D18").Value = "29.298.77"
$ws.Range("E18").Value = "  -0.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "232.15"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("E20").Value = "  -0.28%  "

$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.340"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9997"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.08"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.472"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("E26").Value = "  -2.08%  "

$ws.Range("E27").Value = "  -0.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07350"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +14.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.452"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.475"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.037"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.034"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.817"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.139"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6960"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.571"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01833"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("E38").Value = "  -0.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.848"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.61%  "

$ws.Range("D40").Value = "1.233.43"
$ws.Range("E40").Value = "  -2.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9369"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9998"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").Value = "1.988.01"
$ws.Range("E43").Value = "  -0.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.63"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.96%  "

$ws.Range("E45").Value = "  -1.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000119"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.702"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.945"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1140"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.858"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3896"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.21%  "
